$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.660.82"
$ws.Range("E2").Value = "  +3.81%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.277.58"
$ws.Range("E3").Value = "  +2.36%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.20"
$ws.Range("E5").Value = "  +0.89%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.30"
$ws.Range("E6").Value = "  +5.70%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.595"
$ws.Range("E7").Value = "  +0.42%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.574"
$ws.Range("E9").Value = "  +1.78%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.65"
$ws.Range("E10").Value = "  +3.29%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +1.25%  "

# Row 12
$ws.Range("E12").Value = "  +0.98%  "

# Row 13
$ws.Range("E13").Value = "  +0.59%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.886"
$ws.Range("E14").Value = "  +2.20%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.624.94"
$ws.Range("E15").Value = "  +2.31%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.65"
$ws.Range("E16").Value = "  +2.38%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.280.14"
$ws.Range("E17").Value = "  +2.72%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.460.10"
$ws.Range("E18").Value = "  +3.58%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.06"
$ws.Range("E19").Value = "  -6.37%  "

# Row 20
$ws.Range("E20").Value = "  +4.43%  "

# Row 21
$ws.Range("E21").Value = "  +1.45%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.64"
$ws.Range("E22").Value = "  +1.46%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.21"
$ws.Range("E23").Value = "  +1.10%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.61"
$ws.Range("E24").Value = "  +1.19%  "

# Row 25
$ws.Range("E25").Value = "  +2.69%  "

# Row 26
$ws.Range("E26").Value = "  -0.25%  "

# Row 27
$ws.Range("E27").Value = "  +1.46%  "

# Row 28
$ws.Range("E28").Value = "  +0.29%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.16"
$ws.Range("E29").Value = "  +10.64%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.53"
$ws.Range("E30").Value = "  +2.48%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.21"
$ws.Range("E31").Value = "  +4.54%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.73"
$ws.Range("E32").Value = "  +0.61%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0889"
$ws.Range("E33").Value = "  -1.39%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.77"
$ws.Range("E34").Value = "  -0.48%  "

# Row 35
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.118"
$ws.Range("E35").Value = "  +12.18%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.04"
$ws.Range("E36").Value = "  +3.22%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.17"
$ws.Range("E37").Value = "  -0.45%  "

# Row 38
$ws.Range("E38").Value = "  +0.03%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.94"
$ws.Range("E39").Value = "  +1.16%  "

# Row 40
$ws.Range("E40").Value = "  -0.41%  "

# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0330"
$ws.Range("E41").Value = "  +0.70%  "

# Row 42
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.55"
$ws.Range("E42").Value = "  +23.26%  "

# Row 43
$ws.Range("E43").Value = "  +0.13%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.781.28"
$ws.Range("E44").Value = "  -7.75%  "

# Row 45
$ws.Range("E45").Value = "  -0.25%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.15"
$ws.Range("E46").Value = "  -2.46%  "

# Row 47
$ws.Range("E47").Value = "  +1.71%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "60.45"
$ws.Range("E48").Value = "  -0.49%  "

# Row 49
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.72"
$ws.Range("E49").Value = "  +5.85%  "

# Row 50
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.06"
$ws.Range("E50").Value = "  -2.84%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.55"
$ws.Range("E51").Value = "  +0.72%  "
